$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-22 Monday", "2024-01-23 Tuesday"),
    @("75×84=", "79×75="),
    @("76×37=", "38×56="),
    @("18×51=", "75×25="),
    @("52×61=", "84×12="),
    @("57×70=", "56×61="),
    @("97×63=", "17×78="),
    @("28×51=", "60×59="),
    @("50×90=", "25×60="),
    @("62×58=", "57×56="),
    @("46×99=", "54×70="),
    @("18×80=", "39×35="),
    @("59×99=", "82×41="),
    @("25×42=", "69×64="),
    @("68×95=", "41×20="),
    @("17×16=", "97×27="),
    @("29×45=", "46×24="),
    @("39×74=", "65×58="),
    @("35×89=", "34×94="),
    @("21×34=", "81×83="),
    @("54×75=", "59×49="),
    @("52×69=", "11×88="),
    @("25×32=", "13×15="),
    @("73×67=", "42×81="),
    @("30×39=", "40×26="),
    @("94×48=", "88×56=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
